# Generate Report for Handback
# Updates the localization-status report: the de-de handback round-trip
# completed, so both locale statuses flip from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamps advance,
# the stale-handback-version error clears, and the Status / Error Detail
# columns are widened/narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: status text for both locales ---
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de status columns so the longer text fits.
$overview.Range("E1").ColumnWidth = 29.144371396019366
$overview.Range("F1").ColumnWidth = 29.144371396019366

# --- zh-cn sheet: per-file Status mirrors the same "Handed back" text,
#     refreshed handback datetime, and the stale-version error clears ---
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-09-05 13:01:46"
$zhcn.Range("P2").Value = ""

$zhcn.Range("C1").ColumnWidth = 29.144371396019366
$zhcn.Range("P1").ColumnWidth = 12.913719540550566

# --- de-de sheet: per-file Status mirrors the same "Handed back" text,
#     refreshed handback datetime, and the stale-version error clears ---
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-09-05 13:01:54"
$dede.Range("P2").Value = ""

$dede.Range("C1").ColumnWidth = 29.144371396019366
$dede.Range("P1").ColumnWidth = 12.913719540550566
